$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.07195966666666666
$ws.Range("H2").Value = 0.215879
$ws.Range("I2").Value = 0.07530091904660251
$ws.Range("J2").Value = 0.07530091904660252
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01935066666666667
$ws.Range("N2").Value = 0.058052
$ws.Range("O2").Value = 0.0002219742535102441
$ws.Range("P2").Value = 0.0002219742535102442
$ws.Range("Q2").Value = 0.001392467523111111
$ws.Range("R2").Value = 0.012532207708
$ws.Range("S2").Value = 0.00001671486529400491
$ws.Range("T2").Value = 0.00001671486529400492

# Row 3
$ws.Range("G3").Value = 0.07195966666666666
$ws.Range("H3").Value = 0.215879
$ws.Range("I3").Value = 0.07530091904660251
$ws.Range("J3").Value = 0.07530091904660252
$ws.Range("M3").Value = 87.15592466666668
$ws.Range("N3").Value = 261.467774
$ws.Range("O3").Value = 0.9997780257464898
$ws.Range("P3").Value = 0.9997780257464899
$ws.Range("Q3").Value = 6.271711287038444
$ws.Range("R3").Value = 56.445401583346
$ws.Range("S3").Value = 0.0752842041813085
$ws.Range("T3").Value = 0.07528420418130853

# Row 4
$ws.Range("I4").Value = 0.2743421080169271
$ws.Range("J4").Value = 0.2743421080169271
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01935066666666667
$ws.Range("N4").Value = 0.058052
$ws.Range("O4").Value = 0.0002219742535102441
$ws.Range("P4").Value = 0.0002219742535102442
$ws.Range("Q4").Value = 0.005073144929333332
$ws.Range("R4").Value = 0.04565830436399999
$ws.Range("S4").Value = 0.00006089688463348415
$ws.Range("T4").Value = 0.00006089688463348416

# Row 5
$ws.Range("I5").Value = 0.2743421080169271
$ws.Range("J5").Value = 0.2743421080169271
$ws.Range("M5").Value = 87.15592466666668
$ws.Range("N5").Value = 261.467774
$ws.Range("O5").Value = 0.9997780257464898
$ws.Range("P5").Value = 0.9997780257464899
$ws.Range("Q5").Value = 22.84958161393534
$ws.Range("R5").Value = 205.646234525418
$ws.Range("S5").Value = 0.2742812111322936
$ws.Range("T5").Value = 0.2742812111322936

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6214993333333333
$ws.Range("H6").Value = 1.864498
$ws.Range("I6").Value = 0.6503569729364704
$ws.Range("J6").Value = 0.6503569729364704
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01935066666666667
$ws.Range("N6").Value = 0.058052
$ws.Range("O6").Value = 0.0002219742535102441
$ws.Range("P6").Value = 0.0002219742535102442
$ws.Range("Q6").Value = 0.01202642643288889
$ws.Range("R6").Value = 0.108237837896
$ws.Range("S6").Value = 0.000144362503582755
$ws.Range("T6").Value = 0.0001443625035827551

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6214993333333333
$ws.Range("H7").Value = 1.864498
$ws.Range("I7").Value = 0.6503569729364704
$ws.Range("J7").Value = 0.6503569729364704
$ws.Range("M7").Value = 87.15592466666668
$ws.Range("N7").Value = 261.467774
$ws.Range("O7").Value = 0.9997780257464898
$ws.Range("P7").Value = 0.9997780257464899
$ws.Range("Q7").Value = 54.16734907638356
$ws.Range("R7").Value = 487.506141687452
$ws.Range("S7").Value = 0.6502126104328877
$ws.Range("T7").Value = 0.6502126104328878
